$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column from 2023-09-13 (45182) to 2023-09-15 (45184)
# for all data rows (C2:C46)
$ws.Range("C2:C46").Value = 45184
